# Auto-generated Excel COM-interop script
# Applies the scheduled-runner market-data refresh to the Kujata_Profits workbook.
# For each affected leve row, updates columns H-N (current market price /
# profit calculations) to the freshly fetched values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 937.5
$ws.Range("I17").Value = 500
$ws.Range("K17").Value = 1500
$ws.Range("M17").Value = -1332

$ws.Range("H112").Value = 2310.2354
$ws.Range("J112").Value = 2626.7856
$ws.Range("L112").Value = 7880.3568
$ws.Range("N112").Value = -10096.3568

$ws.Range("H138").Value = 2332.963
$ws.Range("I138").Value = 3479.1667
$ws.Range("J138").Value = 2241.2666
$ws.Range("K138").Value = 10437.5001
$ws.Range("L138").Value = 6723.7998
$ws.Range("M138").Value = -5297.500100000001
$ws.Range("N138").Value = -17003.7998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1340.6923
$ws.Range("I2").Value = 1029.875
$ws.Range("J2").Value = 1838
$ws.Range("K2").Value = 1029.875
$ws.Range("L2").Value = 1838
$ws.Range("M2").Value = -916.875
$ws.Range("N2").Value = -2064

$ws.Range("H25").Value = 1744.3334
$ws.Range("I25").Value = 1899.875
$ws.Range("J25").Value = 500
$ws.Range("K25").Value = 1899.875
$ws.Range("L25").Value = 500
$ws.Range("M25").Value = -1497.875
$ws.Range("N25").Value = -1304

$ws.Range("H32").Value = 6395.977
$ws.Range("I32").Value = 6461.0464
$ws.Range("K32").Value = 6461.0464
$ws.Range("M32").Value = -6174.0464

$ws.Range("H74").Value = 3098.6365
$ws.Range("I74").Value = 2183.5715
$ws.Range("K74").Value = 2183.5715
$ws.Range("M74").Value = -1309.5715

$ws.Range("H77").Value = 3098.6365
$ws.Range("I77").Value = 2183.5715
$ws.Range("K77").Value = 10917.8575
$ws.Range("M77").Value = -6549.8575

$ws.Range("H116").Value = 1340.6923
$ws.Range("I116").Value = 1029.875
$ws.Range("J116").Value = 1838
$ws.Range("K116").Value = 1029.875
$ws.Range("L116").Value = 1838
$ws.Range("M116").Value = 1264.125
$ws.Range("N116").Value = -6426

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1340.6923
$ws.Range("I3").Value = 1029.875
$ws.Range("J3").Value = 1838
$ws.Range("K3").Value = 1029.875
$ws.Range("L3").Value = 1838
$ws.Range("M3").Value = -915.875
$ws.Range("N3").Value = -2066

$ws.Range("H20").Value = 1814.0526
$ws.Range("I20").Value = 1948.6428
$ws.Range("K20").Value = 1948.6428
$ws.Range("M20").Value = -1701.6428

$ws.Range("H26").Value = 7500
$ws.Range("I26").Value = 7500
$ws.Range("K26").Value = 7500
$ws.Range("M26").Value = -7208

$ws.Range("H37").Value = 2575
$ws.Range("I37").Value = 100
$ws.Range("K37").Value = 100
$ws.Range("M37").Value = 37

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 83334696
$ws.Range("I16").Value = 125001500
$ws.Range("K16").Value = 125001500
$ws.Range("M16").Value = -125001213

$ws.Range("H31").Value = 1204.014
$ws.Range("I31").Value = 1084.9193
$ws.Range("J31").Value = 2024.4445
$ws.Range("K31").Value = 1084.9193
$ws.Range("L31").Value = 2024.4445
$ws.Range("M31").Value = -789.9193
$ws.Range("N31").Value = -2614.4445

$ws.Range("H34").Value = 1204.014
$ws.Range("I34").Value = 1084.9193
$ws.Range("J34").Value = 2024.4445
$ws.Range("K34").Value = 1084.9193
$ws.Range("L34").Value = 2024.4445
$ws.Range("M34").Value = -882.9193
$ws.Range("N34").Value = -2428.4445

$ws.Range("H99").Value = 1804.0588
$ws.Range("I99").Value = 1772.0714
$ws.Range("K99").Value = 1772.0714
$ws.Range("M99").Value = -274.0714

$ws.Range("H103").Value = 16814.75
$ws.Range("I103").Value = 12336.333
$ws.Range("J103").Value = 30250
$ws.Range("K103").Value = 12336.333
$ws.Range("L103").Value = 30250
$ws.Range("M103").Value = -11164.333
$ws.Range("N103").Value = -32594

$ws.Range("H113").Value = 83334696
$ws.Range("I113").Value = 125001500
$ws.Range("K113").Value = 125001500
$ws.Range("M113").Value = -124999330

$ws.Range("H126").Value = 1804.0588
$ws.Range("I126").Value = 1772.0714
$ws.Range("K126").Value = 5316.2142
$ws.Range("M126").Value = -2846.2142

$ws.Range("H132").Value = 1705.0625
$ws.Range("I132").Value = 1259.3043
$ws.Range("K132").Value = 3777.9129
$ws.Range("M132").Value = -1247.9129

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1027.1818
$ws.Range("J5").Value = 915.7143
$ws.Range("L5").Value = 2747.1429
$ws.Range("N5").Value = -2971.1429

$ws.Range("H63").Value = 10949.4375
$ws.Range("I63").Value = 2844.5715
$ws.Range("K63").Value = 8533.7145
$ws.Range("M63").Value = -7784.7145

$ws.Range("H66").Value = 10949.4375
$ws.Range("I66").Value = 2844.5715
$ws.Range("K66").Value = 25601.1435
$ws.Range("M66").Value = -21857.1435

$ws.Range("H113").Value = 679.67645
$ws.Range("I113").Value = 510.33334
$ws.Range("K113").Value = 1531.00002
$ws.Range("M113").Value = 638.9999800000001

$ws.Range("H131").Value = 38520016
$ws.Range("J131").Value = 100999.93
$ws.Range("L131").Value = 302999.79
$ws.Range("N131").Value = -313079.79

$ws.Range("H135").Value = 1027.1818
$ws.Range("J135").Value = 915.7143
$ws.Range("L135").Value = 8241.4287
$ws.Range("N135").Value = -13311.4287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H102").Value = 1569.6666
$ws.Range("I102").Value = 1526.0435
$ws.Range("J102").Value = 1713
$ws.Range("K102").Value = 1526.0435
$ws.Range("L102").Value = 1713
$ws.Range("M102").Value = 95.95650000000001
$ws.Range("N102").Value = -4957

$ws.Range("H113").Value = 1086.2941
$ws.Range("I113").Value = 887.8
$ws.Range("J113").Value = 1369.8572
$ws.Range("K113").Value = 887.8
$ws.Range("L113").Value = 1369.8572
$ws.Range("M113").Value = 1282.2
$ws.Range("N113").Value = -5709.8572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3948.8333
$ws.Range("I40").Value = 2655.2856
$ws.Range("K40").Value = 2655.2856
$ws.Range("M40").Value = -2519.2856

$ws.Range("H61").Value = 1192.3636
$ws.Range("I61").Value = 1038.5
$ws.Range("K61").Value = 1038.5
$ws.Range("M61").Value = -836.5

$ws.Range("H113").Value = 1192.3636
$ws.Range("I113").Value = 1038.5
$ws.Range("K113").Value = 1038.5
$ws.Range("M113").Value = 1131.5

$ws.Range("H122").Value = 31251688
$ws.Range("I122").Value = 50001500
$ws.Range("K122").Value = 150004500
$ws.Range("M122").Value = -150002050

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 511.05264
$ws.Range("I113").Value = 367.75
$ws.Range("J113").Value = 756.7143
$ws.Range("K113").Value = 1103.25
$ws.Range("L113").Value = 2270.1429
$ws.Range("M113").Value = 1066.75
$ws.Range("N113").Value = -6610.1429

$ws.Range("H136").Value = 1259.8334
$ws.Range("I136").Value = 867.41174
$ws.Range("K136").Value = 2602.23522
$ws.Range("M136").Value = -52.23522000000003

